$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns retain text formatting (values look numeric e.g. "306.78")
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '44.257.71'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').Value = '2.249.23'
$ws.Range('E3').Value = '  +0.81%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '306.78'
$ws.Range('E5').Value = '  -2.29%  '
$ws.Range('D6').Value = '96.28'
$ws.Range('E6').Value = '  -2.10%  '
$ws.Range('D7').Value = '0.575'
$ws.Range('E7').Value = '  +0.93%  '
$ws.Range('E8').Value = '  +0.21%  '
$ws.Range('D9').Value = '0.529'
$ws.Range('E9').Value = '  -0.95%  '
$ws.Range('D10').Value = '35.38'
$ws.Range('E10').Value = '  -2.15%  '
$ws.Range('D11').Value = '0.0819'
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').Value = '7.29'
$ws.Range('E12').Value = '  -1.41%  '
$ws.Range('D13').Value = '0.105'
$ws.Range('E13').Value = '  +0.27%  '
$ws.Range('D14').Value = '2.592.25'
$ws.Range('E14').Value = '  +0.80%  '
$ws.Range('D15').Value = '2.318.52'
$ws.Range('E15').Value = '  +2.30%  '
$ws.Range('D16').Value = '0.839'
$ws.Range('E16').Value = '  -0.01%  '
$ws.Range('D17').Value = '13.70'
$ws.Range('E17').Value = '  -2.94%  '
$ws.Range('D18').Value = '44.126.28'
$ws.Range('E18').Value = '  +1.09%  '
$ws.Range('D19').Value = '0.0₃0978'
$ws.Range('E19').Value = '  +1.29%  '
$ws.Range('D20').Value = '12.28'
$ws.Range('E20').Value = '  -4.86%  '
$ws.Range('D21').Value = '6.42'
$ws.Range('E21').Value = '  +1.51%  '
$ws.Range('D22').Value = '65.68'
$ws.Range('E22').Value = '  +0.71%  '
$ws.Range('B23').Value = 'PancakeSwap'
$ws.Range('C23').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D23').Value = '3.13'
$ws.Range('E23').Value = '  +4.74%  '
$ws.Range('B24').Value = 'BitcoinCash'
$ws.Range('C24').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D24').Value = '237.63'
$ws.Range('E24').Value = '  +1.51%  '
$ws.Range('D25').Value = '2.01'
$ws.Range('E25').Value = '  -1.28%  '
$ws.Range('E26').Value = '  +0.07%  '
$ws.Range('D27').Value = '10.02'
$ws.Range('E27').Value = '  -0.26%  '
$ws.Range('D28').Value = '2.21'
$ws.Range('E28').Value = '  +1.65%  '
$ws.Range('D29').Value = '38.06'
$ws.Range('E29').Value = '  +4.44%  '
$ws.Range('D30').Value = '6.02'
$ws.Range('E30').Value = '  +1.07%  '
$ws.Range('D31').Value = '20.19'
$ws.Range('E31').Value = '  +1.29%  '
$ws.Range('D32').Value = '153.09'
$ws.Range('E32').Value = '  -3.71%  '
$ws.Range('E33').Value = '  -3.22%  '
$ws.Range('D34').Value = '3.31'
$ws.Range('E34').Value = '  +3.41%  '
$ws.Range('D35').Value = '2.62'
$ws.Range('E35').Value = '  -2.66%  '
$ws.Range('B36').Value = 'Stellar'
$ws.Range('C36').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D36').Value = '0.120'
$ws.Range('E36').Value = '  +2.69%  '
$ws.Range('B37').Value = 'Kaspa'
$ws.Range('C37').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D37').Value = '0.109'
$ws.Range('E37').Value = '  -0.34%  '
$ws.Range('D38').Value = '1.78'
$ws.Range('E38').Value = '  -6.64%  '
$ws.Range('D39').Value = '3.47'
$ws.Range('E39').Value = '  -3.26%  '
$ws.Range('D40').Value = '3.90'
$ws.Range('E40').Value = '  -2.91%  '
$ws.Range('D41').Value = '14.66'
$ws.Range('E41').Value = '  -7.10%  '
$ws.Range('D42').Value = '0.0301'
$ws.Range('E42').Value = '  -2.09%  '
$ws.Range('D43').Value = '1.01'
$ws.Range('E43').Value = '  +0.21%  '
$ws.Range('D44').Value = '1.743.62'
$ws.Range('E44').Value = '  +1.56%  '
$ws.Range('D45').Value = '83.93'
$ws.Range('E45').Value = '  +2.25%  '
$ws.Range('D46').Value = '0.193'
$ws.Range('E46').Value = '  -1.06%  '
$ws.Range('D47').Value = '100.33'
$ws.Range('E47').Value = '  -1.92%  '
$ws.Range('D48').Value = '4.95'
$ws.Range('E48').Value = '  -3.33%  '
$ws.Range('D49').Value = '8.18'
$ws.Range('E49').Value = '  +1.83%  '
$ws.Range('D50').Value = '55.10'
$ws.Range('E50').Value = '  -2.47%  '
$ws.Range('E51').Value = '  -6.52%  '
